# ADD results from server
# Replace the per-year investment-cost breakdown (row 1 headers + row 2 values)
# on every "year" worksheet with the new server results. The "gt"/"dgt"
# (gas turbine / diesel generator turbine) columns are dropped and two new
# technology columns - "gb" (gas boiler) and "btes" (borehole thermal energy
# storage) - are introduced.

$wb = $excel.ActiveWorkbook

$headers = @("eb","gb","hp","st","wi","ieh","chp","ac","ab_ct","ab_hp","cp_ct","cp_hp","ttes","btes","ites")

# Row-2 values (one array per worksheet, in tab order) for columns A..O
# matching $headers above.
$valuesBySheet = @{
    1 = @(39063.99109145206, 0, 483537.6274462014, 0, 2897240.114301849, 94331.34471502228, 0, 25342.77928792104, 0, 0, 0, 0, 0, 23638.06126801545, 19940.13531829346)
    2 = @(30846.52922536713, 0, 1495599.874611417, 0, 0, 70193.79982138964, 0, 56602.42752520426, 0, 0, 0, 0, 0, 51649.16401227913, 42574.77934331147)
    3 = @(242452.4252219552, 0, 943335.270081223, 0, 0, 1425.925979620855, 0, 39373.98526588717, 0, 0, 0, 0, 0, 53308.16490721726, 30023.09380555204)
    4 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 11578.49752443177, 0)
    5 = @(76705.58894163162, 1930.947398408091, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 28147.3462746636, 8312.661449003012)
    6 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
}

for ($s = 1; $s -le $wb.Worksheets.Count; $s++) {
    $ws = $wb.Worksheets.Item($s)

    for ($col = 1; $col -le $headers.Length; $col++) {
        $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
    }

    $rowValues = $valuesBySheet[$s]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item(2, $col).Value = $rowValues[$col - 1]
    }
}
